$d = $word.ActiveDocument

# The 3rd paragraph reads "ตาราง … Sequence Diagram" and is built out of
# the runs (character offsets below are for the *original* document):
#   [33-38)  "ตาราง"    (w:hint="cs", w:cs/)
#   [38-39)  " "         (w:hint="cs", w:cs/)        -> becomes "ที่ " (hint dropped)
#   [39-41)  "… "         (no hint)                   -> splits into "1" + " "
#   [41-49)  "Sequence"   (no hint)
#   [49-57)  " Diagram"   (no hint)
#
# This runtime coalesces adjacent runs whose formatting ends up identical
# whenever a run is touched, so edits are ordered/shaped to land on the
# exact run layout the target diff expects.

# --- 1) " " (right after "ตาราง") -> "ที่ ", dropping w:hint="cs" ---
# Re-assigning the (unchanged) font name makes the engine rewrite the run's
# rFonts without the stale w:hint="cs" marker while keeping w:cs/.
$rHeadSpace = $d.Range(38, 39)
$rHeadSpace.Font.Name = "TH Sarabun New"
$rHeadSpace.Text = "ที่ "

# Text above grew from 1 to 4 characters, so everything from here on is
# shifted by +3 relative to the original offsets.

# --- 2) protect the right-hand boundary ("Sequence") so our edit below
#        does not get coalesced into it ---
$rSeq = $d.Range(44, 52)
$rSeq.Font.Size = 17

# --- 3) "… " -> "1 " ---
$rDots = $d.Range(42, 44)
$rDots.Text = "1 "

# --- 4) split "1 " into separate "1" and " " runs ---
$rSpace = $d.Range(43, 44)
$rSpace.Font.Size = 18
$rSpace.Text = " "
$rSpace2 = $d.Range(43, 44)
$rSpace2.Font.Size = 16

# --- 5) restore "Sequence" formatting ---
$rSeq2 = $d.Range(44, 52)
$rSeq2.Font.Size = 16

Write-Output "done"
